$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect so the cells can be edited.
$ws.Unprotect()

# Update the confidential date note text (shared string)
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-25 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-7
$ws.Range("D2").Value = 0.2477138639934497
$ws.Range("E2").Value = -0.001103576510498083

$ws.Range("D3").Value = 0.4978186270506204
$ws.Range("E3").Value = 0.01101072840203265

$ws.Range("D4").Value = 0.09805727931128262
$ws.Range("E4").Value = 0.005391953546246331

$ws.Range("D5").Value = 0.09972959587334713
$ws.Range("E5").Value = 0.01425601425601419

$ws.Range("D6").Value = 0.05668063377130032
$ws.Range("E6").Value = 0.02394195888754536

$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 0.008515486733206679

# Restore sheet protection as it was before the edit (same flags as the
# original: contents/objects/scenarios protected, column & row formatting
# still allowed). The original password is a legacy hash we cannot recover,
# so the sheet is re-protected without a password.
$ws.Protect($null, $true, $true, $true, $false, $false, $true, $true, $false, $false, $false, $false, $false, $false, $false, $false)
